$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the organism/program/classification values in column B (rows 1-6),
# leaving only the row labels in column A.
$ws.Range("B1").ClearContents()
$ws.Range("B2").ClearContents()
$ws.Range("B3").ClearContents()
$ws.Range("B4").ClearContents()
$ws.Range("B5").ClearContents()
$ws.Range("B6").ClearContents()

# Update the "Fecha de elaboración" timestamp.
$ws.Range("B10").Value = "07-07-2022 11:42:47 am"

# Insert a new "Componente" row above the existing "Fin" row (row 12),
# pushing the "Fin" row down to row 13.
$ws.Rows.Item(12).Insert()

# The inserted row picks up formatting from the row above (the header row);
# reset it back to the default/normal style used by the rest of the data rows.
$ws.Range("A12:G12").Style = "Normal"

$ws.Range("A12").Value = "Componente"
$ws.Range("B12").Value = 5612
$ws.Range("C12").Value = "6 Gobierno Ciudadano"
$ws.Range("D12").Value = "Sumatoria de POA alineados con el Plan Estatal de Desarrollo 2021-2027"
$ws.Range("E12").ClearContents()
$ws.Range("F12").Value = "aa"
$ws.Range("G12").Value = "Secretaría de Planeación y Particiación Ciudadana"
